$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Row 11 (Marking): B11 4 -> 5, C11 -1 -> -1.2
$ws.Range("B11").Value = 5
$ws.Range("C11").Value = -1.2

# Row 12 (Total): B12 76 -> 95, C12 -4 -> -4.8, E12 "72/112" -> "90.2/140"
$ws.Range("B12").Value = 95
$ws.Range("C12").Value = -4.8
$ws.Range("E12").Value = "90.2/140"
